# Apply the edits described by the diff:
#  - Row 2 (existing "ريد بل - 250 مل"): Price D2 changes from 1040 to 1065
#  - New Row 3: Product 5152 "ريد بل فرى شوجر - 250 مل", Packing Unit 2, Price 1065, Visibility YES
#  - New Row 4: Product 7630 "فيورى جولد - 400 مل", Packing Unit 2, Price 205, Visibility YES

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 price
$ws.Range("D2").Value2 = 1065

# New row 3
$ws.Range("A3").Value2 = 5152
$ws.Range("B3").Value2 = "ريد بل فرى شوجر - 250 مل"
$ws.Range("C3").Value2 = 2
$ws.Range("D3").Value2 = 1065

# New row 4
$ws.Range("A4").Value2 = 7630
$ws.Range("B4").Value2 = "فيورى جولد - 400 مل"
$ws.Range("C4").Value2 = 2
$ws.Range("D4").Value2 = 205

# Visibility column (set after the new product names so the shared-string
# table picks up the new names ahead of re-used values)
$ws.Range("E2").Value2 = "YES"
$ws.Range("E3").Value2 = "YES"
$ws.Range("E4").Value2 = "YES"
